# ============================================================================
# Edits A&Dguide.docx per the commit:
#  1. Remove the stray _GoBack bookmark left near the top of the "Game
#     Overview" intro paragraph.
#  2. In the "Demon Lord / High Angel" stat block, move the
#     "Attack Range:<tab>3" paragraph so it appears *before* the
#     "Special:" paragraph (it previously sat between "Special:" and
#     "Movement:"), and leave a fresh _GoBack bookmark at the end of the
#     relocated paragraph (this is where Word's last-edit bookmark ends
#     up after such a move).
#  3. Under "Special Types" -> "Splash:", re-type the trailing sentence
#     ("Attacks effect a cross 2 area from the initial point.") as three
#     separate runs, and add a new "Projectile:" paragraph describing the
#     Projectile special type right after it.
#  4. Drop the stale <w:lastRenderedPageBreak/> marker on the "Game Grid"
#     heading.
# ============================================================================

$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Remove the old _GoBack bookmark.
# ----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ----------------------------------------------------------------------
# 2) Relocate "Attack Range:<tab>3" ahead of "Special:...Splash, Projectile"
#    within the Demon Lord / High Angel block, re-adding _GoBack at the
#    tail of the moved paragraph.
# ----------------------------------------------------------------------
$paras = $d.Paragraphs
$specialIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -eq ("Special:" + [char]9 + [char]9 + "Splash, Projectile " + [char]13)) {
        $specialIdx = $i
        break
    }
}

if ($specialIdx -gt 0) {
    $pSpecial = $paras.Item($specialIdx)
    $pAttackRange = $paras.Item($specialIdx + 1)

    if ($pAttackRange.Range.Text -eq ("Attack Range:" + [char]9 + "3" + [char]13)) {
        $moveRange = $d.Range($pSpecial.Range.Start, $pAttackRange.Range.End)
        $moveXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
                   '<w:r><w:t>Attack Range:</w:t></w:r>' +
                   '<w:r><w:tab/><w:t>3</w:t></w:r>' +
                   '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
                   '<w:bookmarkEnd w:id="0"/>' +
                   '</w:p>' +
                   '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
                   '<w:r><w:t>Special:</w:t></w:r>' +
                   '<w:r><w:tab/></w:r>' +
                   '<w:r><w:tab/></w:r>' +
                   '<w:r><w:t>Splash</w:t></w:r>' +
                   '<w:r><w:t>, Projectile</w:t></w:r>' +
                   '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
                   '</w:p>'
        $moveRange.InsertXML($moveXml)
    }
}

# ----------------------------------------------------------------------
# 3) Re-split the Splash sentence into 3 runs and insert the new
#    "Projectile:" paragraph right after it.
# ----------------------------------------------------------------------
$paras = $d.Paragraphs
$splashIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -eq ("Splash:" + [char]9 + [char]9 + "Attacks affect a cross 2 area from the initial point." + [char]13)) {
        $splashIdx = $i
        break
    }
}

if ($splashIdx -gt 0) {
    $pSplash = $paras.Item($splashIdx)
    $splashXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                 '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
                 '<w:r><w:t>Splash:</w:t></w:r>' +
                 '<w:r><w:tab/></w:r>' +
                 '<w:r><w:tab/></w:r>' +
                 '<w:r><w:t xml:space="preserve">Attacks </w:t></w:r>' +
                 '<w:r><w:t>e</w:t></w:r>' +
                 '<w:r><w:t>ffect a cross 2 area from the initial point.</w:t></w:r>' +
                 '</w:p>' +
                 '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                 '<w:pPr><w:spacing w:after="0"/></w:pPr>' +
                 '<w:r><w:t>Projectile:</w:t></w:r>' +
                 '<w:r><w:tab/><w:t>Can target any unit in attack range (no collision).</w:t></w:r>' +
                 '</w:p>'
    $pSplash.Range.InsertXML($splashXml)
}

# ----------------------------------------------------------------------
# 4) Drop the lastRenderedPageBreak marker on the "Game Grid" heading.
# ----------------------------------------------------------------------
$paras = $d.Paragraphs
$gridIdx = 0
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -eq ("Game Grid" + [char]13)) {
        $gridIdx = $i
        break
    }
}

if ($gridIdx -gt 0) {
    $pGrid = $paras.Item($gridIdx)
    $gridXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' +
               '<w:r><w:t>Game Grid</w:t></w:r>' +
               '</w:p>'
    $pGrid.Range.InsertXML($gridXml)
}
